# 工作進度.xlsx — "Add files via upload" edit
# Adds three new status rows for 黃容 (Huang Rong):
#   Row 11 — under the 12/2 week section (between row 10 and row 13)
#   Row 14 — under the 12/10 week section (row 14 was a blank spacer row)
#   Row 20 — a brand-new 12/23 week section appended after the last row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 11: continuation of the 12/2(一） weekly block (rows 8-10)
# ---------------------------------------------------------------------
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("A11").Value = "12/2(一）"
$ws.Range("B11").Value = "黃容"
$ws.Range("C11").Value = "負責第一部分整理數據"
$ws.Range("E11").Value = "1. 將數據遺漏值填補、計算投資報酬，以及將投資報酬標準化。"
$ws.Range("E11").Font.Name = "Microsoft JhengHei"
$ws.Range("E11").Font.Size = 12
$ws.Range("E11").Font.Color = 0

$ws.Range("A11").RowHeight = 17

# ---------------------------------------------------------------------
# Row 14: continuation of the 12/10(一） weekly block (previously a
# blank spacer row holding only a formatted, empty E14 cell)
# ---------------------------------------------------------------------
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("A14").Value = "12/10(一）"
$ws.Range("B14").Value = "黃容"
$ws.Range("C14").Value = "mySQL, DB使用方法整理"
$ws.Range("E14").Value = "上傳DB以及mysql語法教學"

# ---------------------------------------------------------------------
# Row 20: new 12/23(二) weekly block appended after the last row (19)
# ---------------------------------------------------------------------
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("B19").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("E19").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("A20").Value = "12/23(二)"
$ws.Range("B20").Value = "黃容"
$ws.Range("C20").Value = "書面報告完成"
$ws.Range("E20").Value = "書面報告完成、整理教學。"

# ---------------------------------------------------------------------
# Restore the view state: selection on A14 with the frozen header pane
# scrolled back to the top (matches the author's last on-screen state)
# ---------------------------------------------------------------------
$ws.Range("A14").Select()
